# Arreglo de bugs en excepciones Diamante Bruto y Diamante Corte
$wb = $excel.ActiveWorkbook

# --- Sheet "Sheet": add clients in rows 13 and 14 ---
$wsClientes = $wb.Worksheets.Item("Sheet")

$wsClientes.Range("A13").Value = "Santiago Arango"
$wsClientes.Range("B13").Value = 12
$wsClientes.Range("C13").Value = "santi@gmail.com"
$wsClientes.Range("D13").Value = "Medellín"
$wsClientes.Range("E13").Value = 7155934

$wsClientes.Range("A14").Value = "Bayron Valdés "
$wsClientes.Range("B14").Value = 13
$wsClientes.Range("C14").Value = "bayron2813@gmail.com"
$wsClientes.Range("D14").Value = "Medellín"
$wsClientes.Range("E14").Value = 3015168866

# --- Sheet "Pedidos": add orders in rows 8 through 12 ---
$wsPedidos = $wb.Worksheets.Item("Pedidos")

$wsPedidos.Range("A8").Value = 12
$wsPedidos.Range("B8").Value = 7
$wsPedidos.Range("C8").Value = "El tamaño del Diamante es 0.3 `nEl grabado del Diamante es True `nEl origen del diamante es cabello `nEl tamaño del Diamante es 0.5 `nEl grabado del Diamante es True `nEl origen del diamante es cenizas `nEl corte del diamante es corazon `n"
$wsPedidos.Range("D8").Value = "31/12/2022"

$wsPedidos.Range("A9").Value = 13
$wsPedidos.Range("B9").Value = 8
$wsPedidos.Range("C9").Value = "El tamaño del Diamante es 0.8 `nEl grabado del Diamante es True `nEl origen del diamante es cenizas `nEl tamaño del Diamante es 0.5 `nEl grabado del Diamante es True `nEl origen del diamante es cenizas `nEl corte del diamante es Corazon `n"
$wsPedidos.Range("D9").Value = "23/04/2021"

$wsPedidos.Range("A10").Value = 13
$wsPedidos.Range("B10").Value = 9
$wsPedidos.Range("C10").Value = "El tamaño del Diamante es 0.6 `nEl grabado del Diamante es True `nEl origen del diamante es cenizas `nEl tamaño del Diamante es 0.5 `nEl grabado del Diamante es True `nEl origen del diamante es cenizas `nEl corte del diamante es corazon `n"
$wsPedidos.Range("D10").Value = "20/04/2022"

$wsPedidos.Range("A11").Value = 10
$wsPedidos.Range("B11").Value = 10
$wsPedidos.Range("C11").Value = "El tamaño del Diamante es 0.5 `nEl grabado del Diamante es True `nEl origen del diamante es cenizas `n"
$wsPedidos.Range("D11").Value = "20/04/2022"

$wsPedidos.Range("A12").Value = 10
$wsPedidos.Range("B12").Value = 11
$wsPedidos.Range("C12").Value = "El tamaño del Diamante es 0.6 `nEl grabado del Diamante es True `nEl origen del diamante es Cabello `n"
$wsPedidos.Range("D12").Value = "21/04/2022"
